# edit.ps1
# Applies: adds a new "time_taken" column (F) with per-row timestamps,
# and inserts a new gene row (VPS50) at row 97, shifting the remaining
# gene rows down by one (confidence values follow the shifted rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "time_taken" header in F1, matching the style of the other headers ---
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# --- Populate F2:F116 with the recorded time_taken timestamps ---
$timeTaken = @(
    "2021-10-05 13:40:55.557435",
    "2021-10-05 13:40:55.557446",
    "2021-10-05 13:40:55.557449",
    "2021-10-05 13:40:55.557452",
    "2021-10-05 13:40:55.557455",
    "2021-10-05 13:40:55.557458",
    "2021-10-05 13:40:55.557460",
    "2021-10-05 13:40:55.557463",
    "2021-10-05 13:40:55.557465",
    "2021-10-05 13:40:55.557468",
    "2021-10-05 13:40:55.557470",
    "2021-10-05 13:40:55.557473",
    "2021-10-05 13:40:55.557476",
    "2021-10-05 13:40:55.557478",
    "2021-10-05 13:40:55.557481",
    "2021-10-05 13:40:55.557483",
    "2021-10-05 13:40:55.557486",
    "2021-10-05 13:40:55.557489",
    "2021-10-05 13:40:55.557491",
    "2021-10-05 13:40:55.557494",
    "2021-10-05 13:40:55.557496",
    "2021-10-05 13:40:55.557499",
    "2021-10-05 13:40:55.557501",
    "2021-10-05 13:40:55.557504",
    "2021-10-05 13:40:55.557507",
    "2021-10-05 13:40:55.557509",
    "2021-10-05 13:40:55.557512",
    "2021-10-05 13:40:55.557514",
    "2021-10-05 13:40:55.557517",
    "2021-10-05 13:40:55.557519",
    "2021-10-05 13:40:55.557522",
    "2021-10-05 13:40:55.557524",
    "2021-10-05 13:40:55.557527",
    "2021-10-05 13:40:55.557530",
    "2021-10-05 13:40:55.557532",
    "2021-10-05 13:40:55.557535",
    "2021-10-05 13:40:55.557537",
    "2021-10-05 13:40:55.557540",
    "2021-10-05 13:40:55.557542",
    "2021-10-05 13:40:55.557545",
    "2021-10-05 13:40:55.557548",
    "2021-10-05 13:40:55.557551",
    "2021-10-05 13:40:55.557553",
    "2021-10-05 13:40:55.557556",
    "2021-10-05 13:40:55.557558",
    "2021-10-05 13:40:55.557560",
    "2021-10-05 13:40:55.557563",
    "2021-10-05 13:40:55.557565",
    "2021-10-05 13:40:55.557568",
    "2021-10-05 13:40:55.557570",
    "2021-10-05 13:40:55.557573",
    "2021-10-05 13:40:55.557575",
    "2021-10-05 13:40:55.557578",
    "2021-10-05 13:40:55.557581",
    "2021-10-05 13:40:55.557583",
    "2021-10-05 13:40:55.557586",
    "2021-10-05 13:40:55.557588",
    "2021-10-05 13:40:55.557591",
    "2021-10-05 13:40:55.557593",
    "2021-10-05 13:40:55.557595",
    "2021-10-05 13:40:55.557598",
    "2021-10-05 13:40:55.557600",
    "2021-10-05 13:40:55.557603",
    "2021-10-05 13:40:55.557605",
    "2021-10-05 13:40:55.557609",
    "2021-10-05 13:40:55.557611",
    "2021-10-05 13:40:55.557614",
    "2021-10-05 13:40:55.557616",
    "2021-10-05 13:40:55.557619",
    "2021-10-05 13:40:55.557621",
    "2021-10-05 13:40:55.557624",
    "2021-10-05 13:40:55.557626",
    "2021-10-05 13:40:55.557629",
    "2021-10-05 13:40:55.557631",
    "2021-10-05 13:40:55.557634",
    "2021-10-05 13:40:55.557636",
    "2021-10-05 13:40:55.557640",
    "2021-10-05 13:40:55.557643",
    "2021-10-05 13:40:55.557646",
    "2021-10-05 13:40:55.557648",
    "2021-10-05 13:40:55.557651",
    "2021-10-05 13:40:55.557653",
    "2021-10-05 13:40:55.557656",
    "2021-10-05 13:40:55.557658",
    "2021-10-05 13:40:55.557661",
    "2021-10-05 13:40:55.557663",
    "2021-10-05 13:40:55.557666",
    "2021-10-05 13:40:55.557668",
    "2021-10-05 13:40:55.557671",
    "2021-10-05 13:40:55.557673",
    "2021-10-05 13:40:55.557676",
    "2021-10-05 13:40:55.557678",
    "2021-10-05 13:40:55.557682",
    "2021-10-05 13:40:55.557685",
    "2021-10-05 13:40:55.557687",
    "2021-10-05 13:40:55.557690",
    "2021-10-05 13:40:55.557692",
    "2021-10-05 13:40:55.557695",
    "2021-10-05 13:40:55.557697",
    "2021-10-05 13:40:55.557700",
    "2021-10-05 13:40:55.557703",
    "2021-10-05 13:40:55.557705",
    "2021-10-05 13:40:55.557707",
    "2021-10-05 13:40:55.557710",
    "2021-10-05 13:40:55.557712",
    "2021-10-05 13:40:55.557715",
    "2021-10-05 13:40:55.557717",
    "2021-10-05 13:40:55.557720",
    "2021-10-05 13:40:55.557724",
    "2021-10-05 13:40:55.557727",
    "2021-10-05 13:40:55.557730",
    "2021-10-05 13:40:55.557732",
    "2021-10-05 13:40:55.557735",
    "2021-10-05 13:40:55.557737",
    "2021-10-05 13:40:55.557740"
)
for ($i = 0; $i -lt $timeTaken.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

# --- A new gene (VPS50) was inserted at row 97, pushing rows 97-115 down to 98-116. ---
# --- Target gene symbol / name / confidence for rows 97-116 (columns B, C, D) ---
$geneRows = @(
    @{ Row = 97; B = "VPS50"; C = "VPS50, EARP/GARPII complex subunit"; D = "2" },
    @{ Row = 98; B = "EOMES"; C = "eomesodermin"; D = "1" },
    @{ Row = 99; B = "ERMARD"; C = "ER membrane associated RNA degradation"; D = "1" },
    @{ Row = 100; B = "GMPPB"; C = "GDP-mannose pyrophosphorylase B"; D = "1" },
    @{ Row = 101; B = "MCF2"; C = "MCF.2 cell line derived transforming sequence"; D = "1" },
    @{ Row = 102; B = "PEX11A"; C = "peroxisomal biogenesis factor 11 alpha"; D = "1" },
    @{ Row = 103; B = "POMK"; C = "protein-O-mannose kinase"; D = "1" },
    @{ Row = 104; B = "TSC2"; C = "TSC complex subunit 2"; D = "1" },
    @{ Row = 105; B = "ATP1A3"; C = "ATPase Na+/K+ transporting subunit alpha 3"; D = "0" },
    @{ Row = 106; B = "ENO1"; C = "enolase 1"; D = "0" },
    @{ Row = 107; B = "GRIN2B"; C = "glutamate ionotropic receptor NMDA type subunit 2B"; D = "0" },
    @{ Row = 108; B = "MAPK8IP3"; C = "mitogen-activated protein kinase 8 interacting protein 3"; D = "0" },
    @{ Row = 109; B = "NPRL2"; C = "NPR2 like, GATOR1 complex subunit"; D = "0" },
    @{ Row = 110; B = "NPRL3"; C = "NPR3 like, GATOR1 complex subunit"; D = "0" },
    @{ Row = 111; B = "PTEN"; C = "phosphatase and tensin homolog"; D = "0" },
    @{ Row = 112; B = "RAB18"; C = "RAB18, member RAS oncogene family"; D = "0" },
    @{ Row = 113; B = "RAB3GAP1"; C = "RAB3 GTPase activating protein catalytic subunit 1"; D = "0" },
    @{ Row = 114; B = "RAB3GAP2"; C = "RAB3 GTPase activating non-catalytic protein subunit 2"; D = "0" },
    @{ Row = 115; B = "SCN3A"; C = "sodium voltage-gated channel alpha subunit 3"; D = "0" },
    @{ Row = 116; B = "SLC35A2"; C = "solute carrier family 35 member A2"; D = "0" }
)
foreach ($entry in $geneRows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    # Column D holds small integers stored as *text*; prefix with a quote char so
    # Excel keeps it as a string instead of inferring a Number, then restore the
    # plain "Normal" style so no stray number-format / quote-prefix marker sticks.
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = "'" + $entry.D
    $dCell.Style = "Normal"
}

# --- Row 116 is brand new: copy formatting from row 115, then set its values ---
$ws.Range("A115:E115").Copy() | Out-Null
$ws.Range("A116:E116").PasteSpecial(-4122) | Out-Null
$ws.Range("A116").Value = 114
$ws.Range("E116").Value = "Malformations of cortical development"

